$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.464.45'
$ws.Range("E2").Value = '  +1.81%  '

$ws.Range("D3").Value = '1.668.40'
$ws.Range("E3").Value = '  +1.31%  '

$ws.Range("D4").Value = "'0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = "'237.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.70%  '

$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("D7").Value = "'0.4788"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.91%  '

$ws.Range("D8").Value = "'0.2625"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.99%  '

$ws.Range("D9").Value = "'0.06174"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.82%  '

$ws.Range("B10").Value = 'WrappedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D10").Value = '1.669.03'
$ws.Range("E10").Value = '  +1.26%  '

$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = "'0.06977"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.04%  '

$ws.Range("D12").Value = "'14.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.23%  '

$ws.Range("D13").Value = "'0.5873"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.32%  '

$ws.Range("D14").Value = "'4.368"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.24%  '

$ws.Range("D15").Value = "'74.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.76%  '

$ws.Range("E16").Value = '  -0.06%  '

$ws.Range("D17").Value = "'0.9998"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.05%  '

$ws.Range("D18").Value = '25.465.68'
$ws.Range("E18").Value = '  +1.84%  '

$ws.Range("D19").Value = "'0.000006746"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.89%  '

$ws.Range("D20").Value = "'11.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.36%  '

$ws.Range("D21").Value = '1.881.75'
$ws.Range("E21").Value = '  +1.40%  '

$ws.Range("D22").Value = "'4.444"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.49%  '

$ws.Range("D23").Value = "'8.729"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.19%  '

$ws.Range("D24").Value = "'5.280"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.27%  '

$ws.Range("D25").Value = "'137.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.02%  '

$ws.Range("D26").Value = "'15.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.03%  '

$ws.Range("E27").Value = '  -0.73%  '

$ws.Range("D28").Value = "'1.718"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.84%  '

$ws.Range("D29").Value = "'104.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.90%  '

$ws.Range("D30").Value = "'3.937"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.75%  '

$ws.Range("D31").Value = "'0.07786"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.82%  '

$ws.Range("D32").Value = "'3.640"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.65%  '

$ws.Range("E33").Value = '  +0.00%  '

$ws.Range("D34").Value = "'0.04220"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.07%  '

$ws.Range("D35").Value = "'2.603"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.34%  '

$ws.Range("D36").Value = "'0.6092"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.62%  '

$ws.Range("D37").Value = "'0.9479"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.72%  '

$ws.Range("D38").Value = "'2.598"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.84%  '

$ws.Range("D39").Value = "'0.8552"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.83%  '

$ws.Range("D40").Value = "'1.000"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.09%  '

$ws.Range("D41").Value = "'1.850"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.67%  '

$ws.Range("D42").Value = "'0.01476"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.89%  '

$ws.Range("D43").Value = "'96.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.08%  '

$ws.Range("D44").Value = "'0.3768"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.41%  '

$ws.Range("D45").Value = "'4.832"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.08%  '

$ws.Range("D46").Value = "'0.1118"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.76%  '

$ws.Range("D47").Value = "'6.183"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.28%  '

$ws.Range("D48").Value = "'0.05251"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.19%  '

$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = "'7.380"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.94%  '
